$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2341
$ws.Range("I28").Value = 2250.5
$ws.Range("J28").Value = 2612.5
$ws.Range("K28").Value = 2250.5
$ws.Range("L28").Value = 2612.5
$ws.Range("M28").Value = -1765.5
$ws.Range("N28").Value = -3582.5
$ws.Range("H43").Value = 257526.06
$ws.Range("I43").Value = 1062.8572
$ws.Range("K43").Value = 1062.8572
$ws.Range("M43").Value = -993.8571999999999
$ws.Range("H64").Value = 7497.8
$ws.Range("J64").Value = 7496.3335
$ws.Range("L64").Value = 7496.3335
$ws.Range("N64").Value = -7992.3335
$ws.Range("H67").Value = 7497.8
$ws.Range("J67").Value = 7496.3335
$ws.Range("L67").Value = 7496.3335
$ws.Range("N67").Value = -9212.333500000001
$ws.Range("H74").Value = 83343496
$ws.Range("J74").Value = 20252
$ws.Range("L74").Value = 20252
$ws.Range("N74").Value = -22124
$ws.Range("H76").Value = 4664.1665
$ws.Range("I76").Value = 4248
$ws.Range("K76").Value = 4248
$ws.Range("M76").Value = -3933
$ws.Range("H77").Value = 83343496
$ws.Range("J77").Value = 20252
$ws.Range("L77").Value = 101260
$ws.Range("N77").Value = -110620
$ws.Range("H79").Value = 4664.1665
$ws.Range("I79").Value = 4248
$ws.Range("K79").Value = 4248
$ws.Range("M79").Value = -3156
$ws.Range("H80").Value = 30050.53
$ws.Range("I80").Value = 11690.889
$ws.Range("K80").Value = 35072.667
$ws.Range("M80").Value = -34074.667
$ws.Range("H83").Value = 30050.53
$ws.Range("I83").Value = 11690.889
$ws.Range("K83").Value = 105218.001
$ws.Range("M83").Value = -100226.001
$ws.Range("H86").Value = 56220376
$ws.Range("I86").Value = 75002960
$ws.Range("K86").Value = 75002960
$ws.Range("M86").Value = -75001837
$ws.Range("H89").Value = 56220376
$ws.Range("I89").Value = 75002960
$ws.Range("K89").Value = 375014800
$ws.Range("M89").Value = -375009184
$ws.Range("H103").Value = 1241.2632
$ws.Range("J103").Value = 1331.25
$ws.Range("L103").Value = 3993.75
$ws.Range("N103").Value = -5165.75
$ws.Range("H106").Value = 978.61536
$ws.Range("I106").Value = 978.61536
$ws.Range("K106").Value = 978.61536
$ws.Range("M106").Value = -347.61536
$ws.Range("H132").Value = 1459
$ws.Range("I132").Value = 1489.9354
$ws.Range("K132").Value = 4469.8062
$ws.Range("M132").Value = -1939.8062
$ws.Range("H137").Value = 2713.0435
$ws.Range("I137").Value = 2570
$ws.Range("K137").Value = 7710
$ws.Range("M137").Value = -5160
$ws.Range("H138").Value = 3965.4167
$ws.Range("I138").Value = 1110.425
$ws.Range("K138").Value = 3331.275
$ws.Range("M138").Value = 1808.725
$ws.Range("H141").Value = 7938219.5
$ws.Range("I141").Value = 8334930.5
$ws.Range("K141").Value = 25004791.5
$ws.Range("M141").Value = -24999611.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 33531.668
$ws.Range("I5").Value = 33531.668
$ws.Range("K5").Value = 33531.668
$ws.Range("M5").Value = -33419.668
$ws.Range("I32").Value = 1739980
$ws.Range("J32").Value = 6916
$ws.Range("K32").Value = 1739980
$ws.Range("L32").Value = 6916
$ws.Range("M32").Value = -1739693
$ws.Range("N32").Value = -7490
$ws.Range("H74").Value = 17720.87
$ws.Range("I74").Value = 22595.357
$ws.Range("J74").Value = 4072.3
$ws.Range("K74").Value = 22595.357
$ws.Range("L74").Value = 4072.3
$ws.Range("M74").Value = -21721.357
$ws.Range("N74").Value = -5820.3
$ws.Range("H77").Value = 17720.87
$ws.Range("I77").Value = 22595.357
$ws.Range("J77").Value = 4072.3
$ws.Range("K77").Value = 112976.785
$ws.Range("L77").Value = 20361.5
$ws.Range("M77").Value = -108608.785
$ws.Range("N77").Value = -29097.5
$ws.Range("H132").Value = 4996.709
$ws.Range("I132").Value = 3492.9143
$ws.Range("J132").Value = 7628.35
$ws.Range("K132").Value = 10478.7429
$ws.Range("L132").Value = 22885.05
$ws.Range("M132").Value = -7948.742899999999
$ws.Range("N132").Value = -27945.05

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 33531.668
$ws.Range("I4").Value = 33531.668
$ws.Range("K4").Value = 33531.668
$ws.Range("M4").Value = -33416.668
$ws.Range("H22").Value = 270.42856
$ws.Range("I22").Value = 269
$ws.Range("K22").Value = 269
$ws.Range("M22").Value = -96
$ws.Range("H86").Value = 38504040
$ws.Range("I86").Value = 62398.883
$ws.Range("K86").Value = 62398.883
$ws.Range("M86").Value = -61275.883
$ws.Range("H89").Value = 38504040
$ws.Range("I89").Value = 62398.883
$ws.Range("K89").Value = 311994.415
$ws.Range("M89").Value = -306378.415
$ws.Range("H99").Value = 3249523.2
$ws.Range("I99").Value = 2706.7778
$ws.Range("K99").Value = 2706.7778
$ws.Range("M99").Value = -1208.7778
$ws.Range("H134").Value = 4637.2856
$ws.Range("I134").Value = 1316.2258
$ws.Range("K134").Value = 3948.6774
$ws.Range("M134").Value = -1413.6774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8831.270500000001
$ws.Range("I31").Value = 4022.3572
$ws.Range("K31").Value = 4022.3572
$ws.Range("M31").Value = -3727.3572
$ws.Range("H34").Value = 8831.270500000001
$ws.Range("I34").Value = 4022.3572
$ws.Range("K34").Value = 4022.3572
$ws.Range("M34").Value = -3820.3572
$ws.Range("H99").Value = 9268.538
$ws.Range("I99").Value = 10296.167
$ws.Range("K99").Value = 10296.167
$ws.Range("M99").Value = -8798.166999999999
$ws.Range("H126").Value = 9268.538
$ws.Range("I126").Value = 10296.167
$ws.Range("K126").Value = 30888.501
$ws.Range("M126").Value = -28418.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7031.467
$ws.Range("J34").Value = 8743.583000000001
$ws.Range("L34").Value = 26230.749
$ws.Range("N34").Value = -26398.749
$ws.Range("H39").Value = 12313.883
$ws.Range("J39").Value = 13899.066
$ws.Range("L39").Value = 41697.198
$ws.Range("N39").Value = -42285.198
$ws.Range("H55").Value = 68258560
$ws.Range("J55").Value = 8341019.5
$ws.Range("L55").Value = 25023058.5
$ws.Range("N55").Value = -25023412.5
$ws.Range("H122").Value = 1769212.2
$ws.Range("I122").Value = 4715588
$ws.Range("J122").Value = 1386.7
$ws.Range("K122").Value = 42440292
$ws.Range("L122").Value = 12480.3
$ws.Range("M122").Value = -42437842
$ws.Range("N122").Value = -17380.3
$ws.Range("H131").Value = 1543.125
$ws.Range("I131").Value = 1466.5
$ws.Range("J131").Value = 1589.1
$ws.Range("K131").Value = 4399.5
$ws.Range("L131").Value = 4767.299999999999
$ws.Range("M131").Value = 640.5
$ws.Range("N131").Value = -14847.3
$ws.Range("H137").Value = 2784.111
$ws.Range("I137").Value = 1699.3334
$ws.Range("J137").Value = 3326.5
$ws.Range("K137").Value = 5098.0002
$ws.Range("L137").Value = 9979.5
$ws.Range("M137").Value = 1.999799999999595
$ws.Range("N137").Value = -20179.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10327.6
$ws.Range("I70").Value = 8864.375
$ws.Range("K70").Value = 8864.375
$ws.Range("M70").Value = -8594.375
$ws.Range("H73").Value = 10327.6
$ws.Range("I73").Value = 8864.375
$ws.Range("K73").Value = 8864.375
$ws.Range("M73").Value = -7928.375
$ws.Range("H80").Value = 5746
$ws.Range("I80").Value = 3619
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 3619
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -2621
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 5746
$ws.Range("I83").Value = 3619
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 18095
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -13103
$ws.Range("N83").Value = -59984
$ws.Range("H113").Value = 6604.8887
$ws.Range("I113").Value = 3224.0833
$ws.Range("J113").Value = 8295.291999999999
$ws.Range("K113").Value = 3224.0833
$ws.Range("L113").Value = 8295.291999999999
$ws.Range("M113").Value = -1054.0833
$ws.Range("N113").Value = -12635.292
$ws.Range("H122").Value = 5511723.5
$ws.Range("I122").Value = 7959712
$ws.Range("K122").Value = 23879136
$ws.Range("M122").Value = -23876686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 938.4167
$ws.Range("I16").Value = 906.44446
$ws.Range("K16").Value = 906.44446
$ws.Range("M16").Value = -736.44446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1448.4445
$ws.Range("I100").Value = 1157.6
$ws.Range("K100").Value = 2315.2
$ws.Range("M100").Value = -1774.2
$ws.Range("H132").Value = 9812575
$ws.Range("I132").Value = 12198589
$ws.Range("J132").Value = 29919.2
$ws.Range("K132").Value = 36595767
$ws.Range("L132").Value = 89757.60000000001
$ws.Range("M132").Value = -36593237
$ws.Range("N132").Value = -94817.60000000001
$ws.Range("H136").Value = 25030170
$ws.Range("I136").Value = 58824316
$ws.Range("J136").Value = 51890.086
$ws.Range("K136").Value = 176472948
$ws.Range("L136").Value = 155670.258
$ws.Range("M136").Value = -176470398
